# Regenerate the "K" column (column G) values for rows 2-21 of Sheet1.
# The underlying save_data was regenerated to use K (strike count proxy)
# instead of the old Strike# values, so the numeric values in column G
# are replaced with their newly computed counterparts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 2
    3  = 1
    4  = 0
    5  = 3
    6  = 4
    7  = 0
    8  = 2
    9  = 0
    10 = 1
    11 = 0
    12 = 1
    13 = 4
    14 = 5
    15 = 1
    16 = 2
    17 = 1
    18 = 3
    19 = 0
    20 = 1
    21 = 4
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
